# Commit: fixed a typo "0.35" -> "35.0"
#
# Two edits from the source diff:
#   1. Slide 7 title: the three runs that spell out
#      "Evaluation (ACL 2012" + ") " + "\u2013 SVM" collapse into a single
#      run with the same overall text ("Evaluation (ACL 2012) \u2013 SVM").
#      Re-selecting the whole run (via Characters) and re-assigning the
#      same text forces PowerPoint to coalesce the runs.
#   2. Slide 8 table ("Table 4"): the "All Positive" row's Precision cell
#      had a typo, "0.35", which should read "35.0".

$p = $ppt.ActivePresentation

# --- Slide 7: merge the title's runs ----------------------------------------
$s7 = $p.Slides.Item(7)
$titleRange = $s7.Shapes.Item(1).TextFrame.TextRange
$titleText = $titleRange.Text
$titleRange.Characters(1, $titleText.Length).Text = $titleText

# --- Slide 8: fix the "0.35" -> "35.0" typo in the results table -----------
$s8 = $p.Slides.Item(8)
$table = $s8.Shapes.Item(3).Table
$table.Cell(4, 2).Shape.TextFrame.TextRange.Text = "35.0"
